# Refreshed crypto price/volume(1h) figures scraped on 2023-11-13.
# A handful of rows also changed which coin they list (the source feed
# re-ranked a few coins), so B (Coin) and C (Link) are rewritten too where
# that happened.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of new Price values are plain numbers ("1.00", "0.999", ...).
# Assigning those bare to Range.Value lets Excel coerce them to real numbers
# and silently eat the significant trailing/leading zeros, so they are
# written with a leading apostrophe (quote-prefix) to force literal text,
# exactly as if a user had typed them into a Text-looking cell.
$ws.Range('D2').Value = '36.842.16'
$ws.Range('E2').Value = '  -0.98%  '

$ws.Range('D3').Value = '2.090.88'
$ws.Range('E3').Value = '  +1.43%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.19%  '

$ws.Range('D5').Value = '''246.60'
$ws.Range('E5').Value = '  -1.09%  '

$ws.Range('D6').Value = '''0.654'
$ws.Range('E6').Value = '  -2.01%  '

$ws.Range('D7').Value = '''0.998'
$ws.Range('E7').Value = '  -0.16%  '

$ws.Range('D8').Value = '''56.63'
$ws.Range('E8').Value = '  -4.86%  '

$ws.Range('D9').Value = '''59.15'
$ws.Range('E9').Value = '  -2.22%  '

$ws.Range('D10').Value = '''0.370'
$ws.Range('E10').Value = '  -4.35%  '

$ws.Range('D11').Value = '''0.0773'
$ws.Range('E11').Value = '  -2.05%  '

$ws.Range('E12').Value = '  +0.93%  '

$ws.Range('D13').Value = '''15.17'
$ws.Range('E13').Value = '  -4.76%  '

$ws.Range('D14').Value = '''0.888'
$ws.Range('E14').Value = '  +6.23%  '

$ws.Range('D15').Value = '2.379.02'
$ws.Range('E15').Value = '  +1.03%  '

$ws.Range('D16').Value = '''5.61'
$ws.Range('E16').Value = '  -3.27%  '

$ws.Range('D17').Value = '2.102.78'
$ws.Range('E17').Value = '  +2.38%  '

$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').Value = '''17.63'
$ws.Range('E18').Value = '  -2.95%  '

$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '36.831.98'
$ws.Range('E19').Value = '  -1.10%  '

$ws.Range('D20').Value = '''73.41'
$ws.Range('E20').Value = '  -2.34%  '

$ws.Range('D21').Value = '0.0₃0882'
$ws.Range('E21').Value = '  -2.44%  '

$ws.Range('D22').Value = '''5.50'
$ws.Range('E22').Value = '  +1.36%  '

$ws.Range('D23').Value = '''236.81'
$ws.Range('E23').Value = '  -0.67%  '

$ws.Range('D24').Value = '''1.00'
$ws.Range('E24').Value = '  +0.17%  '

$ws.Range('D25').Value = '''2.43'
$ws.Range('E25').Value = '  -2.97%  '

$ws.Range('D26').Value = '''9.96'
$ws.Range('E26').Value = '  +5.54%  '

$ws.Range('D27').Value = '''2.18'
$ws.Range('E27').Value = '  -0.64%  '

$ws.Range('D28').Value = '''167.86'
$ws.Range('E28').Value = '  -0.90%  '

$ws.Range('D29').Value = '''21.21'
$ws.Range('E29').Value = '  +5.06%  '

$ws.Range('D30').Value = '''5.36'
$ws.Range('E30').Value = '  +10.34%  '

$ws.Range('D31').Value = '''0.124'
$ws.Range('E31').Value = '  -1.18%  '

$ws.Range('D32').Value = '''1.20'
$ws.Range('E32').Value = '  +7.07%  '

$ws.Range('D33').Value = '''4.77'
$ws.Range('E33').Value = '  +4.03%  '

$ws.Range('D34').Value = '''0.0613'
$ws.Range('E34').Value = '  -1.86%  '

$ws.Range('D35').Value = '''2.44'
$ws.Range('E35').Value = '  +5.48%  '

$ws.Range('D36').Value = '''0.999'
$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('E37').Value = '  +5.32%  '

$ws.Range('D38').Value = '''0.0851'
$ws.Range('E38').Value = '  -6.36%  '

$ws.Range('D39').Value = '''1.30'
$ws.Range('E39').Value = '  -3.12%  '

$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '''1.17'
$ws.Range('E40').Value = '  +1.29%  '

$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').Value = '''4.95'
$ws.Range('E41').Value = '  -3.61%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.0222'
$ws.Range('E42').Value = '  -0.70%  '

$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').Value = '''2.94'
$ws.Range('E43').Value = '  -5.78%  '

$ws.Range('D44').Value = '''0.0957'
$ws.Range('E44').Value = '  -10.14%  '

$ws.Range('D45').Value = '''97.33'
$ws.Range('E45').Value = '  +0.40%  '

$ws.Range('D46').Value = '''16.50'
$ws.Range('E46').Value = '  -6.09%  '

$ws.Range('D47').Value = '1.345.83'
$ws.Range('E47').Value = '  +4.26%  '

$ws.Range('D48').Value = '''2.41'
$ws.Range('E48').Value = '  -2.97%  '

$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '''7.10'
$ws.Range('E49').Value = '  +2.58%  '

$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').Value = '''2.87'
$ws.Range('E50').Value = '  -1.11%  '

$ws.Range('D51').Value = '2.258.97'
$ws.Range('E51').Value = '  +0.60%  '
